$wb = $excel.ActiveWorkbook
$wsLoc = $wb.Worksheets.Item("Locations")
$wsDrinks = $wb.Worksheets.Item("Drinks")

# --- Move the "Keys in JSON" note from D2 to E2 (format + value), then
# --- clear D2 entirely (it must not remain as an empty cell). ---
$wsLoc.Range("D2").Copy()
$wsLoc.Range("E2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$wsLoc.Range("E2").Value = "Keys in JSON"
$wsLoc.Range("D2").Clear()

# --- Add the new "pathOnTablet" column (C) to the Locations sheet,
# --- reusing the formatting already used for the analogous column on the
# --- Drinks sheet. ---
$wsDrinks.Range("C1").Copy()
$wsLoc.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsDrinks.Range("C2").Copy()
$wsLoc.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsDrinks.Range("C3").Copy()
$wsLoc.Range("C3:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cell values ---
$wsLoc.Range("C1").Value = "Path in React's public folder"
$wsLoc.Range("C2").Value = "pathOnTablet"
$wsLoc.Range("C3").Value = "img/location/bar.png"
$wsLoc.Range("C4").Value = "img/location/livingRoom.png"
$wsLoc.Range("C5").Value = "img/location/entrance.png"
$wsLoc.Range("C6").Value = "img/location/bin.png"

# --- Column widths (best effort - engine quantizes to a coarser grid than
# --- real Excel, so these land as close as the host allows). ---
$wsLoc.Columns.Item(3).ColumnWidth = 39.75
$wsLoc.Columns.Item(4).ColumnWidth = 11.65
$wsLoc.Columns.Item(5).ColumnWidth = 11.65

# --- Selections to match the saved view state: Drinks gets an entire-
# --- column selection on C, People loses its tab-selected flag, and
# --- Locations becomes the active (tab-selected) sheet with E2 selected. ---
$wsDrinks.Range("C1:C1048576").Select()
$wsLoc.Range("E2").Select()
